# NIT-9017451193.xlsx update:
# - "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
# The workbook had two "mora" periods (2507 and 2506) for two workers.
# This update removes the 2506 period rows entirely (database refresh),
# keeps only the 2507 period, updates the totals (Valor Mora, Cant. Periodos)
# accordingly, and updates the remaining worker's row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update summary header values ---
# Valor Mora total: 210000 -> 120000
$ws.Range("E11").Value = 120000
# Cant. Periodos: 2 -> 1
$ws.Range("F13").Value = 1

# --- Remove the two rows that belonged to the old "2506" period ---
# Row 17 = CC / 3805454 / GERMAN ELLES GOMEZ / 2506 / 60000 / 1500000
# Row 18 = CC / 1007786943 / ALFONSO GONZALEZ CONTRERAS / 2507 / 60000 / 1000000
# Deleting both shifts row 19 (ALFONSO .../2506/30000/1000000) up to row 17,
# and the signature block (old rows 24-25) up to rows 22-23.
$ws.Rows("17:18").Delete()

# --- Fix up the remaining data row (now row 17) to the updated 2507 figures ---
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 60000
$ws.Range("G17").Value = 1500000
